# Updated canine keywords; added samples, files obj
# Targets:
#   - CypherOutput sheet (sheet1): rewrite rows 2-5 with refreshed data,
#     add new rows 6-7 (NCATS-COP01-CCB040254 / transcriptomics sample and
#     GLIOMA01-i_4990 / genomics sample).
#   - StatOutput sheet (sheet4): refresh the file/sample/case/study counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CypherOutput ("file name" data export for the American Staffordshire
# Terrier breed filter)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CypherOutput")

# Force column G ("Age") to be written as literal text so values like
# "7.3" / "13.0" / "4.0" keep their exact string form instead of being
# coerced to numbers.
$ageRange = $ws.Range("G2:G7")
$ageRange.NumberFormat = "@"

$ws.Range("A2").Value = "COTC007B-0412"
$ws.Range("B2").Value = "COTC007B"
$ws.Range("C2").Value = "Clinical Trial"
$ws.Range("D2").Value = "American Staffordshire Terrier"
$ws.Range("E2").Value = "Lymphoma"
$ws.Range("F2").Value = "Va"
$ws.Range("G2").Value = "7.3"
$ws.Range("H2").Value = "Male"
$ws.Range("I2").Value = "Yes"

$ws.Range("A3").Value = "COTC007B-0301"
$ws.Range("B3").Value = "COTC007B"
$ws.Range("C3").Value = "Clinical Trial"
$ws.Range("D3").Value = "American Staffordshire Terrier"
$ws.Range("E3").Value = "Lymphoma"
$ws.Range("F3").Value = "IIIa"
$ws.Range("G3").Value = "5.8"
$ws.Range("H3").Value = "Female"
$ws.Range("I3").Value = "Yes"

$ws.Range("A4").Value = "COTC007B-0409"
$ws.Range("B4").Value = "COTC007B"
$ws.Range("C4").Value = "Clinical Trial"
$ws.Range("D4").Value = "American Staffordshire Terrier"
$ws.Range("E4").Value = "Lymphoma"
$ws.Range("F4").Value = "IIIa"
$ws.Range("G4").Value = "2.4"
$ws.Range("H4").Value = "Male"
$ws.Range("I4").Value = "Yes"

$ws.Range("A5").Value = "COTC007B-0617"
$ws.Range("B5").Value = "COTC007B"
$ws.Range("C5").Value = "Clinical Trial"
$ws.Range("D5").Value = "American Staffordshire Terrier"
$ws.Range("E5").Value = "Lymphoma"
$ws.Range("F5").Value = "III"
$ws.Range("G5").Value = "10.0"
$ws.Range("H5").Value = "Female"
$ws.Range("I5").Value = "Yes"

$ws.Range("A6").Value = "NCATS-COP01-CCB040254"
$ws.Range("B6").Value = "NCATS-COP01"
$ws.Range("C6").Value = "Transcriptomics"
$ws.Range("D6").Value = "American Staffordshire Terrier"
$ws.Range("E6").Value = "Pulmonary Neoplasms"
$ws.Range("F6").Value = "Unknown"
$ws.Range("G6").Value = "13.0"
$ws.Range("H6").Value = "Female"
$ws.Range("I6").Value = "Yes"

$ws.Range("A7").Value = "GLIOMA01-i_4990"
$ws.Range("B7").Value = "GLIOMA01"
$ws.Range("C7").Value = "Genomics"
$ws.Range("D7").Value = "American Staffordshire Terrier"
$ws.Range("E7").Value = "Glioma"
$ws.Range("F7").Value = "Unknown"
$ws.Range("G7").Value = "4.0"
$ws.Range("H7").Value = "Male"
$ws.Range("I7").Value = "Yes"

# Drop the temporary text number-format now that the literal strings are
# committed, so no stray cell style is left behind.
$ageRange.ClearFormats()

# ---------------------------------------------------------------------
# StatOutput (file/sample/case/study counts for the same breed filter)
# ---------------------------------------------------------------------
$stat = $wb.Worksheets.Item("StatOutput")

$statRange = $stat.Range("A2:D2")
$statRange.NumberFormat = "@"

$stat.Range("A2").Value = "11"
$stat.Range("B2").Value = "8"
$stat.Range("C2").Value = "6"
$stat.Range("D2").Value = "3"

$statRange.ClearFormats()
